$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new feedback entry was inserted as row 4, pushing the previous
# rows 4-12 down to rows 5-13.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = 45912.44327967593
$ws.Range("D4").Value = "NjAwYjE5OTAtOGViYy00Y2Q1LWI4MzAtNzhkZTRlNjcxOWJiOjU3MDE2"
